$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report formulas & format
# - H1's label stays "TOTAL:" (its shared-string index shifts automatically
#   because the string table is rebuilt once the strings below change).
# - A3 no longer carries the stray "." label.
# - B3 / J3 / N3 / O3 now carry a "date" label (string placeholders for the
#   date-formatted columns under Fecha de pago / checkin / checkout).
# - The sheet's lingering I1 selection is reset back to the default A1 cell.

$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "date"
$ws.Range("J3").Value = "date"
$ws.Range("N3").Value = "date"
$ws.Range("O3").Value = "date"

[void]$ws.Range("A1").Select()
